# Mercury Tours - CamposRegister.xlsx
# Registro de un nuevo usuario de prueba (Fabian Alfonso) en la fila 2.
# "Suite registrar vuelo, falta el cp006."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Fabian"
$ws.Range("B2").Value = "Alfonso"
$ws.Range("C2").Value = 311
$ws.Range("D2").Value = "elkin3001"
$ws.Range("E2").Value = "Galan M#14"
$ws.Range("F2").Value = "Calarca"
$ws.Range("G2").Value = "Quindio"
$ws.Range("H2").Value = 57
$ws.Range("I2").Value = "Colombia"
$ws.Range("J2").Value = "elkin3001"
$ws.Range("K2").Value = "elkin3001"
$ws.Range("L2").Value = "elkin3001"

# Columns whose new content is wider than before get resized (their
# "best fit" auto width is replaced by the resulting fixed width),
# mirroring what Excel does when a column's content forces it to grow.
$ws.Range("D1").ColumnWidth = 7.833333333333333
$ws.Range("E1").ColumnWidth = 8.666666666666666
$ws.Range("F1").ColumnWidth = 6.166666666666667
$ws.Range("I1").ColumnWidth = 7.666666666666667

# After typing the last value of the row, the cursor naturally advances
# to the next cell in the entry flow.
$selected = $ws.Range("L3").Select()

$saved = $wb.Save()
